# Regenerate tables and main figures from final dual-VIIRS dataset
# Update header row labels across all table sheets to human-readable names.

$wb = $excel.ActiveWorkbook

# Table1_country_summary
$ws1 = $wb.Worksheets.Item("Table1_country_summary")
$ws1.Range("A1").Value = "Country (ISO3)"
$ws1.Range("B1").Value = "Total Hansen loss 2015–2023 (Mha)"
$ws1.Range("C1").Value = "Mean annual Hansen loss (kha)"
$ws1.Range("D1").Value = "Mean discrepancy ratio (Hansen/FRA)"
$ws1.Range("E1").Value = "Mean protected-loss share"
$ws1.Range("F1").Value = "Total VIIRS fire detections (millions)"

# Table2_yearly_snapshot
$ws2 = $wb.Worksheets.Item("Table2_yearly_snapshot")
$ws2.Range("A1").Value = "Country (ISO3)"
$ws2.Range("B1").Value = "Year"
$ws2.Range("C1").Value = "Hansen loss (kha)"
$ws2.Range("D1").Value = "FRA net forest area change (kha)"
$ws2.Range("E1").Value = "Discrepancy ratio (Hansen/FRA)"

# Table3_protected_share
$ws3 = $wb.Worksheets.Item("Table3_protected_share")
$ws3.Range("A1").Value = "Country (ISO3)"
$ws3.Range("B1").Value = "Mean protected-loss share"
$ws3.Range("C1").Value = "Maximum protected-loss share"
$ws3.Range("D1").Value = "Year of maximum protected-loss share"

# Table4_fire_loss_corr
$ws4 = $wb.Worksheets.Item("Table4_fire_loss_corr")
$ws4.Range("A1").Value = "Country (ISO3)"
$ws4.Range("B1").Value = "Pearson correlation (Hansen vs VIIRS)"
$ws4.Range("C1").Value = "Slope (ha per detection)"
$ws4.Range("D1").Value = "Intercept (ha)"
$ws4.Range("E1").Value = "p-value"

# Table5_interp_check
$ws5 = $wb.Worksheets.Item("Table5_interp_check")
$ws5.Range("A1").Value = "Country (ISO3)"
$ws5.Range("B1").Value = "Max absolute difference (kha)"
$ws5.Range("C1").Value = "Mean absolute difference (kha)"
